# Adaption of gap_types to 2 groups ("Arbeit" & "Privat") instead of 3 groups
# ("Arbeit", "Soziales Umfeld", "Hausarbeit & Selbstsorge").
#
# The only substantive content change is in worksheet "QAGlist_Teil1", column M
# (Gap1_type): the old 3-way grouping labels
#   "Haushalt & Selbstsorge", "Soziales Umfeld",
#   "Arbeit,Haushalt & Selbstsorge", "Arbeit, Haushalt & Selbstsorge"
# are replaced by the new 2-way grouping labels "Privat" / "Arbeit, Privat".
# Once those old labels are no longer referenced anywhere in the workbook,
# they naturally drop out of the shared-strings table on save (mirroring
# what a human editing this in Excel and re-saving would produce).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QAGlist_Teil1")

# --- Column M ("Gap1_type") relabelling -----------------------------------
# Former "Haushalt & Selbstsorge" / "Soziales Umfeld" rows -> "Privat"
$ws.Range("M3").Value  = "Privat"
$ws.Range("M4").Value  = "Privat"
$ws.Range("M5").Value  = "Privat"
$ws.Range("M6").Value  = "Privat"
$ws.Range("M7").Value  = "Privat"
$ws.Range("M8").Value  = "Privat"
$ws.Range("M17").Value = "Privat"
$ws.Range("M20").Value = "Privat"
$ws.Range("M28").Value = "Privat"
$ws.Range("M29").Value = "Privat"
$ws.Range("M30").Value = "Privat"

# Former "Arbeit,Haushalt & Selbstsorge" / "Arbeit, Haushalt & Selbstsorge"
# rows -> "Arbeit, Privat"
$ws.Range("M21").Value = "Arbeit, Privat"
$ws.Range("M22").Value = "Arbeit, Privat"
$ws.Range("M23").Value = "Arbeit, Privat"
$ws.Range("M26").Value = "Arbeit, Privat"

# --- View state: mirror the saved selection on the sheet -------------------
$ws.Activate()
$ws.Range("K31").Select()
